$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update trip-plan description cells (column C = planned, D = actual)
$ws.Range("C2").Value = "Заезд Milan – Domodosolla -- Bognanco – кэмпинг Yolki Palki"
$ws.Range("D2").Value = "Заезд Milan – Domodosolla -- Bognanco – кэмпинг Yolki Palki"
$ws.Range("C3").Value = "пер. Passo di Campo (2180 н/к) – оз. Lagho di Campo (2290)"
$ws.Range("D3").Value = "пер. Passo di Campo (2180 н/к) – оз. Lagho di Campo (2290)"
$ws.Range("C4").Value = "пер. Passo di Pontimia (н/к, 2387) -- пос. Gmeinalp (1850) – ур. Galki (2300)"
$ws.Range("D4").Value = "пер. Passo di Pontimia (н/к, 2387) -- пос. Gmeinalp (1850) – ур. Galki (2300)"
$ws.Range("C6").Value = "Almagelleralp – дер. Saas-Almagel -- Mischabel camping / Kappelenweg – полуднёвка"
$ws.Range("D6").Value = "Almagelleralp – дер. Saas-Almagel – кемпирг “Am Kappelenweg” – полуднёвка"
$ws.Range("C7").Value = "дер. Saas-Fee – хиж. Mischabelhutte"
$ws.Range("D7").Value = "дер. Saas-Fee – хиж. Mischabelhutte"
$ws.Range("C8").Value = "Хиж. Mischabelhutte -- Windjoch (1Б, 3850, сн-л) -- вер. Ulrichshorn (1Б, 3925) – хиж. Bordierhutte"
$ws.Range("D8").Value = "Хиж. Mischabelhutte -- Windjoch (1Б, 3850, сн-л) -- вер. Ulrichshorn (1Б, 3925) – хиж. Bordierhutte"
$ws.Range("C9").Value = "Хиж. Bordierhutte -- пер. 2664 (н/к) – хиж. Europahutte – кемпинг Tasch"
$ws.Range("D9").Value = "Хиж. Bordierhutte -- пер. 2664 (н/к) – хиж. Europahutte – кемпинг Rand"
$ws.Range("C11").Value = "Zermatt – ст. Furi – ст. Trockener Steg (переезд на подъёмнике) – подход по Breithorn"
$ws.Range("D11").Value = "Zermatt – ст. Furi – ст. Trockener Steg (переезд на подъёмнике *) – пер. P3824"
$ws.Range("C12").Value = "вер. Breithorn (1Б, 4164, снежн.) -- развилка Bivacco Rossi e Volante (3700) -- оз. 2740"
$ws.Range("D12").Value = "вер. Breithorn (1Б, 4164, снежн.) -- развилка Bivacco Rossi e Volante (3700) -- оз. 2740"
$ws.Range("C13").Value = "оз. 2740 – дол. Verraz  – дер. Blanchard -- пер. Сolletto di Nano (1А, 2650)"
$ws.Range("D13").Value = "оз. 2740 – дол. Verraz  – дер. Blanchard -- пер. Сolletto di Nano (1А, 2650)"
$ws.Range("C14").Value = "пер. Col de Nannaz (1А, 2773) – дер. Chamois -- дер. Nuarsaz – кемп. Cervino (1130)"
$ws.Range("C16").Value = "пер. Col de Nannaz (1А, 2773) – пер. Col de Fontaines (н/к, 2696) – дер. Cheneil -- camping Glair"
$ws.Range("D16").Value = "пер. Col de Nannaz (1А, 2773) – вер. Becca Trecare (н/к, 3032) –  пер. Col de Fontaines (н/к, 2696) – дер. Cheneil -- camping Glair"
$ws.Range("C17").Value = "кемп. Glair – дер. Buisson – дер. Grand-Moulin –  гор. St. Vincent"
$ws.Range("D17").Value = "кемп. Glair – дер. Buisson – дер. Grand-Moulin –  переезд в гор. St. Vincent (**)"

# Update numeric distance (column E, L in km) values
$ws.Range("E11").Value = 10.7
$ws.Range("E12").Value = 11.5
$ws.Range("E16").Value = 11.9
$ws.Range("E17").Value = 9.3

# Adjust row heights to fit the revised text
$ws.Rows(2).RowHeight = 23.85
$ws.Rows(4).RowHeight = 23.85
$ws.Rows(6).RowHeight = 35.05
$ws.Rows(8).RowHeight = 35.05
$ws.Rows(11).RowHeight = 35.05
$ws.Rows(14).RowHeight = 35.05
$ws.Rows(17).RowHeight = 23.85

# Restore the active cell selection
$ws.Range("D16").Select() | Out-Null

